# Fruta / hortaliza, semanal
# Weekly refresh of the "Vega Modelo de Temuco - Caigua" price rows.
# The new data reshuffles rows 2-28 (D: Fecha, J: Volumen, K/L/M: Precio
# min/max/promedio, P: Precio $/Kg) while every other column stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current values for the columns that move (D, J, K, L, M, P)
# for every data row (2..28) before writing anything, so source rows are
# not overwritten before they've been read.
$cols = @(4, 10, 11, 12, 13, 16)   # D, J, K, L, M, P
$firstRow = 2
$lastRow = 28

$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# new row number -> old (source) row number
$mapping = @{
    2  = 5
    3  = 16
    4  = 18
    5  = 2
    6  = 22
    7  = 9
    8  = 10
    9  = 26
    10 = 25
    11 = 14
    12 = 11
    13 = 21
    14 = 27
    15 = 8
    16 = 20
    17 = 4
    18 = 17
    19 = 12
    20 = 13
    21 = 6
    22 = 7
    23 = 28
    24 = 3
    25 = 24
    26 = 19
    27 = 15
    28 = 23
}

foreach ($newRow in $mapping.Keys) {
    $srcRow = $mapping[$newRow]
    $src = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($newRow, $c).Value2 = $src[$c]
    }
}
